$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Find rows by account number in column A and apply the edits.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = $lastRow; $r -ge 1; $r--) {
    $acct = $ws.Cells.Item($r, 1).Value2

    if ($acct -eq "004334062") {
        # MERG: Saldo 31000 -> 33855.83
        $ws.Cells.Item($r, 3).Value = 33855.83
    }
    elseif ($acct -eq "004224284") {
        # PRISCILLA row removed entirely
        $ws.Rows.Item($r).Delete()
    }
    elseif ($acct -eq "000806386") {
        # FERNANDA row removed entirely
        $ws.Rows.Item($r).Delete()
    }
}

$wb.Save()
